$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, pushing the old "Totals:" row (and everything
# below it) down by one row. This makes room for a third "Bad Driver" entry.
$ws.Rows("5:5").Insert()

# Row 3: previously "22.190.0.4"; now becomes the "23.40.0.4" entry with
# updated sample counts.
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.40.0.4"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 401
$ws.Range("D3").Value = 96.90000000000001

# Row 4: previously "23.40.0.4"; now becomes the "22.190.0.4" entry with
# updated sample counts.
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.190.0.4"
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 637
$ws.Range("D4").Value = 98

# Row 5: brand-new "Bad Driver" entry.
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 286
$ws.Range("D5").Value = 98.2

# Row 6 (formerly row 5): "Totals:" row, now reflecting all three entries.
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 1324

# Good Drivers table shifted down by the inserted row; update the sample
# counts for the first two drivers (values for the third driver, 22.100.1.1,
# are unchanged).
$ws.Range("B14").Value = 11140
$ws.Range("B15").Value = 14487

# Touch J21 (without altering its appearance) so Excel's used range / sheet
# dimension extends to J21, matching the new overall sheet extent.
$ws.Cells.Item(21, 10).Font.Bold = $false
